$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header style (H1) onto the two new
# header cells, then set their text. Copy/PasteSpecial(formats) reuses the
# existing style index instead of minting a new (near-duplicate) one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I66 and J2:J66
$iVals = @(10,7,7,7,7,8,8,6,7,8,7,7,5,7,5,8,5,6,7,7,5,4,9,5,7,7,7,7,6,6,6,7,9,9,5,4,8,7,5,7,6,9,4,7,5,7,6,6,7,6,7,8,5,6,7,6,8,4,6,7,5,6,6,7,7)
$jVals = @(10,7,7,7,7,8,8,6,8,8,7,7,5,7,5,8,6,6,7,7,5,5,9,5,7,7,7,7,7,6,7,7,9,9,5,4,8,7,5,7,6,9,4,7,5,7,6,6,7,6,7,8,5,6,7,7,8,4,7,7,6,6,6,7,7)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
